$d = $word.ActiveDocument

# Original paragraph text is "Versi" + "on" + " 2" + "." = "Version 2."
# Target paragraph text is "Version" + " 1." = "Version 1."

# 1. Remove the trailing "." run (offsets 9-10), which sits after the bookmark.
$d.Range(9, 10).Delete()

# 2. Change the " 2" run (offsets 7-9) into " 1."
$d.Range(7, 9).Text = " 1."

# 3. Merge "Versi" + "on" into a single "Version" run: delete the "on" run
#    (offsets 5-7) then rewrite the remaining "Versi" run (offsets 0-5) as "Version".
$d.Range(5, 7).Delete()
$d.Range(0, 5).Text = "Version"
